$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "db" (sheet1): add a new "devNonce" header column in H2.
# ---------------------------------------------------------------------
$wsDb = $wb.Worksheets.Item("db")
$wsDb.Range("H2").Value = "devNonce"

# ---------------------------------------------------------------------
# Sheet "knihovny" (sheet2): no textual change, content stays the same
# (only the shared-string table is reshuffled elsewhere).
# ---------------------------------------------------------------------

# ---------------------------------------------------------------------
# Sheet "zdroje" (sheet3): rework into a small key/value table with a
# new "crypto"/"gui" section, drop the now-unused pycryptodome AES link
# and add the new cryptography + DearPyGui references.
# ---------------------------------------------------------------------
$wsZ = $wb.Worksheets.Item("zdroje")

# Drop the hyperlinks that are moving / disappearing before rewriting
# the cells underneath them.
$existingLinks = @($wsZ.Hyperlinks)
foreach ($link in $existingLinks) {
    $addr = $link.Range.Address()
    if ($addr -eq '$B$4' -or $addr -eq '$B$6') {
        $link.Delete()
    }
}

# K2: plain "depricated" note, no hyperlink/style.
$wsZ.Range("K2").Value = "depricated"

# B4 used to be the pycryptodome AES link; it becomes the (now
# unlinked, but still hyperlink-styled) cryptography symmetric
# encryption reference.
$wsZ.Range("B4").Value = "https://cryptography.io/en/latest/hazmat/primitives/symmetric-encryption/"
$wsZ.Range("B4").Style = "Hypertextový odkaz"

# B6 used to hold the fernet link, which now moves to K3; keep the
# hyperlink formatting on the now-empty cell.
$wsZ.Range("B6").Value = ""
$wsZ.Range("B6").Style = "Hypertextový odkaz"

# K3: the fernet link, relocated here with hyperlink styling.
$wsZ.Range("K3").Value = "https://cryptography.io/en/latest/fernet/"
$wsZ.Hyperlinks.Add($wsZ.Range("K3"), "https://cryptography.io/en/latest/fernet/")
$wsZ.Range("K3").Style = "Hypertextový odkaz"

# B5: new DearPyGui reference, with hyperlink.
$wsZ.Range("B5").Value = "https://github.com/hoffstadt/DearPyGui"
$wsZ.Hyperlinks.Add($wsZ.Range("B5"), "https://github.com/hoffstadt/DearPyGui")
$wsZ.Range("B5").Style = "Hypertextový odkaz"

# New label column (A) describing each row.
$wsZ.Range("A2").Value = "sqlite3"
$wsZ.Range("A3").Value = "mikrotik api"
$wsZ.Range("A4").Value = "crypto"
$wsZ.Range("A5").Value = "gui"

# Widen column A to fit its new labels (best-fit-like width).
$wsZ.Columns.Item(1).ColumnWidth = 9.7109375

# Final selection / active sheet: "zdroje" becomes the active tab with
# A6 selected; "db" keeps H2 selected.
$wsDb.Range("H2").Select()
$wsZ.Range("A6").Select()
